$wb = $excel.ActiveWorkbook

# --- Rename the "carNameAndPriceTest" sheet to "carNamesAndCarPricesTest" ---
$wsCars = $wb.Worksheets.Item("carNameAndPriceTest")
$wsCars.Name = "carNamesAndCarPricesTest"

$wsFind = $wb.Worksheets.Item("findNewCars")

# --- Update existing rows: runMode column changes from "N" to "Y" ---
$wsCars.Range("B3").Value = "Y"
$wsCars.Range("B4").Value = "Y"

# --- Append two new rows of test data ---
$wsCars.Range("A5").Value = "chrome"
$wsCars.Range("B5").Value = "Y"
$wsCars.Range("C5").Value = "hyundai"

$wsCars.Range("A6").Value = "chrome"
$wsCars.Range("B6").Value = "Y"
$wsCars.Range("C6").Value = "audi"

# --- Update sheet selections / cursor positions ---
$null = $wsFind.Activate()
$null = $wsFind.Range("D13").Select()

# --- Make carNamesAndCarPricesTest the active/selected sheet ---
$null = $wsCars.Activate()
$null = $wsCars.Range("E7").Select()
